# "Stage" added to data-tab
$wb = $excel.ActiveWorkbook

# --- Sheet "data": insert a new column A ("Stage") before the existing data ---
$wsData = $wb.Worksheets.Item("data")

$wsData.Columns("A:A").Insert()

$wsData.Range("A1").Value = "Stage"
$wsData.Range("A2").Value = "Test"
$wsData.Range("A3").Value = "Test"
$wsData.Range("A4").Value = "Test"

# Hyperlinks are not automatically shifted by the column insert in this
# runtime, so remove and recreate them in their new location (column D).
$wsData.Hyperlinks.Delete()
$wsData.Hyperlinks.Add($wsData.Range("D2"), "https://webdemo.baangt.org/")
$wsData.Hyperlinks.Add($wsData.Range("D3"), "https://webdemo.baangt.org/")
$wsData.Hyperlinks.Add($wsData.Range("D4"), "https://webdemo.baangt.org/")

# Re-apply the original "Link" cell style (Hyperlinks.Add introduces a
# duplicate style record otherwise).
$wsData.Range("D2").Style = "Link"
$wsData.Range("D3").Style = "Link"
$wsData.Range("D4").Style = "Link"

# Update the selection on the "data" sheet without leaving it as the active tab.
$wsData.Range("A5").Select()

# --- Sheet "TestStepExecution": move the selection and restore it as the active tab ---
$wsExec = $wb.Worksheets.Item("TestStepExecution")
$wsExec.Range("D29").Select()
